$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 271.5
$ws.Range("I28").Value = 299.33334
$ws.Range("J28").Value = 188
$ws.Range("K28").Value = 299.33334
$ws.Range("L28").Value = 188
$ws.Range("M28").Value = 185.66666
$ws.Range("N28").Value = -1158
$ws.Range("H40").Value = 5727.273
$ws.Range("I40").Value = 5727.273
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 5727.273
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = -5552.273
$ws.Range("N40").ClearContents()
$ws.Range("H64").Value = 53095.6
$ws.Range("I64").Value = 3102.4443
$ws.Range("J64").Value = 93999.09
$ws.Range("K64").Value = 3102.4443
$ws.Range("L64").Value = 93999.09
$ws.Range("M64").Value = -2854.4443
$ws.Range("N64").Value = -94495.09
$ws.Range("H67").Value = 53095.6
$ws.Range("I67").Value = 3102.4443
$ws.Range("J67").Value = 93999.09
$ws.Range("K67").Value = 3102.4443
$ws.Range("L67").Value = 93999.09
$ws.Range("M67").Value = -2244.4443
$ws.Range("N67").Value = -95715.09
$ws.Range("H76").Value = 3596.6206
$ws.Range("I76").Value = 3526.92
$ws.Range("J76").Value = 4032.25
$ws.Range("K76").Value = 3526.92
$ws.Range("L76").Value = 4032.25
$ws.Range("M76").Value = -3211.92
$ws.Range("N76").Value = -4662.25
$ws.Range("H79").Value = 3596.6206
$ws.Range("I79").Value = 3526.92
$ws.Range("J79").Value = 4032.25
$ws.Range("K79").Value = 3526.92
$ws.Range("L79").Value = 4032.25
$ws.Range("M79").Value = -2434.92
$ws.Range("N79").Value = -6216.25
$ws.Range("H129").Value = 724.4286
$ws.Range("J129").Value = 817.86957
$ws.Range("L129").Value = 2453.60871
$ws.Range("N129").Value = -12453.60871
$ws.Range("H132").Value = 2908067.5
$ws.Range("I132").Value = 3379576.5
$ws.Range("J132").Value = 428.5
$ws.Range("K132").Value = 10138729.5
$ws.Range("L132").Value = 1285.5
$ws.Range("M132").Value = -10136199.5
$ws.Range("N132").Value = -6345.5
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 2928.3333
$ws.Range("I63").Value = 2723.5715
$ws.Range("J63").Value = 3406.111
$ws.Range("K63").Value = 2723.5715
$ws.Range("L63").Value = 3406.111
$ws.Range("M63").Value = -2037.5715
$ws.Range("N63").Value = -4778.111
$ws.Range("H66").Value = 2928.3333
$ws.Range("I66").Value = 2723.5715
$ws.Range("J66").Value = 3406.111
$ws.Range("K66").Value = 13617.8575
$ws.Range("L66").Value = 17030.555
$ws.Range("M66").Value = -10185.8575
$ws.Range("N66").Value = -23894.555
$ws.Range("H88").Value = 5272.5
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 5272.5
$ws.Range("K88").Value = 0
$ws.Range("L88").Value = 5272.5
$ws.Range("M88").ClearContents()
$ws.Range("N88").Value = -6084.5
$ws.Range("H91").Value = 5272.5
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = 5272.5
$ws.Range("K91").Value = 0
$ws.Range("L91").Value = 5272.5
$ws.Range("M91").ClearContents()
$ws.Range("N91").Value = -8080.5
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H19").Value = 17150
$ws.Range("J19").Value = 17150
$ws.Range("L19").Value = 17150
$ws.Range("N19").Value = -17496
$ws.Range("H56").Value = 0
$ws.Range("J56").Value = 0
$ws.Range("L56").Value = 0
$ws.Range("N56").ClearContents()
$ws.Range("H105").Value = 1972.909
$ws.Range("I105").Value = 1942.2222
$ws.Range("K105").Value = 1942.2222
$ws.Range("M105").Value = -195.2221999999999
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2025.1613
$ws.Range("I58").Value = 2127.6
$ws.Range("J58").Value = 1838.909
$ws.Range("K58").Value = 2127.6
$ws.Range("L58").Value = 1838.909
$ws.Range("M58").Value = -1924.6
$ws.Range("N58").Value = -2244.909
$ws.Range("H62").Value = 0
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("M62").ClearContents()
$ws.Range("N62").ClearContents()
$ws.Range("H65").Value = 0
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("M65").ClearContents()
$ws.Range("N65").ClearContents()
$ws.Range("H107").Value = 1087.15
$ws.Range("I107").Value = 1058.3572
$ws.Range("J107").Value = 1154.3334
$ws.Range("K107").Value = 1058.3572
$ws.Range("L107").Value = 1154.3334
$ws.Range("M107").Value = 861.6428000000001
$ws.Range("N107").Value = -4994.3334
$ws.Range("H132").Value = 1412.4857
$ws.Range("I132").Value = 1127.0646
$ws.Range("J132").Value = 3624.5
$ws.Range("K132").Value = 3381.1938
$ws.Range("L132").Value = 10873.5
$ws.Range("M132").Value = -851.1938
$ws.Range("N132").Value = -15933.5
$ws.Range("H136").Value = 2025.1613
$ws.Range("I136").Value = 2127.6
$ws.Range("J136").Value = 1838.909
$ws.Range("K136").Value = 6382.799999999999
$ws.Range("L136").Value = 5516.727000000001
$ws.Range("M136").Value = -3832.799999999999
$ws.Range("N136").Value = -10616.727
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H112").Value = 927346.2
$ws.Range("I112").Value = 509
$ws.Range("J112").Value = 957901.25
$ws.Range("K112").Value = 1527
$ws.Range("L112").Value = 2873703.75
$ws.Range("M112").Value = -419
$ws.Range("N112").Value = -2875919.75
$ws.Range("H122").Value = 1914.4286
$ws.Range("I122").Value = 2142.8572
$ws.Range("J122").Value = 1800.2142
$ws.Range("K122").Value = 19285.7148
$ws.Range("L122").Value = 16201.9278
$ws.Range("M122").Value = -16835.7148
$ws.Range("N122").Value = -21101.9278
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H26").Value = 10000
$ws.Range("J26").Value = 10000
$ws.Range("L26").Value = 10000
$ws.Range("N26").Value = -10560
$ws.Range("H50").Value = 10000
$ws.Range("J50").Value = 10000
$ws.Range("L50").Value = 10000
$ws.Range("N50").Value = -10996
$ws.Range("H70").Value = 4076.6897
$ws.Range("I70").Value = 4080.25
$ws.Range("J70").Value = 4068.7778
$ws.Range("K70").Value = 4080.25
$ws.Range("L70").Value = 4068.7778
$ws.Range("M70").Value = -3810.25
$ws.Range("N70").Value = -4608.7778
$ws.Range("H73").Value = 4076.6897
$ws.Range("I73").Value = 4080.25
$ws.Range("J73").Value = 4068.7778
$ws.Range("K73").Value = 4080.25
$ws.Range("L73").Value = 4068.7778
$ws.Range("M73").Value = -3144.25
$ws.Range("N73").Value = -5940.7778
$ws.Range("H80").Value = 2550.8333
$ws.Range("J80").Value = 2561.2
$ws.Range("L80").Value = 2561.2
$ws.Range("N80").Value = -4557.2
$ws.Range("H83").Value = 2550.8333
$ws.Range("J83").Value = 2561.2
$ws.Range("L83").Value = 12806
$ws.Range("N83").Value = -22790
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H36").Value = 31651.111
$ws.Range("J36").Value = 31651.111
$ws.Range("L36").Value = 31651.111
$ws.Range("N36").Value = -32775.111
$ws.Range("H56").Value = 2683.6667
$ws.Range("I56").Value = 2683.6667
$ws.Range("J56").Value = 0
$ws.Range("K56").Value = 2683.6667
$ws.Range("L56").Value = 0
$ws.Range("M56").Value = -1992.6667
$ws.Range("N56").ClearContents()
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H58").Value = 32300.2
$ws.Range("J58").Value = 32300.2
$ws.Range("L58").Value = 32300.2
$ws.Range("N58").Value = -32916.2
$ws.Range("H61").Value = 5625
$ws.Range("I61").Value = 5625
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 5625
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -5333
$ws.Range("N61").ClearContents()
